# Scrum daily update 10.05.2017
# Update "Effort Actual" (column K) values on the Sprint Backlog sheet
# and move the active selection to K4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint Backlog")
$ws.Activate()

$ws.Range("K4").Value = 3
$ws.Range("K5").Value = 3
$ws.Range("K7").Value = 5
$ws.Range("K8").Value = 4

$ws.Range("K4").Select()
